# Daily AIESEC Global Talent opportunity scrape refresh - 2026-02-12 04:30:11 UTC
# Replaces the 9 data rows (rows 2-10) with freshly scraped opportunities and
# drops the 6 older rows (11-16) that fell out of the scrape window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Drop the trailing rows (11-16) that no longer exist in the new scrape.
#    Do this first so the later writes target a sheet already shaped A1:H10.
# ---------------------------------------------------------------------------
$ws.Rows("11:16").Delete()

# ---------------------------------------------------------------------------
# 2) Move the yellow "PREMIUM = Yes" highlight style from E4 to E2 (the new
#    premium listing is row 2 this time, not row 4). Copy carries the style
#    (and value) over; ClearFormats + re-set the value restores E4 to plain.
# ---------------------------------------------------------------------------
$ws.Range("E4").Copy($ws.Range("E2"))
$ws.Range("E4").ClearFormats()
$ws.Range("E4").Value = "No"

# ---------------------------------------------------------------------------
# 3) Row 2 - [EXP] Digital Marketing PMO Intern (Maastricht, Netherlands)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "'1331871"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = 'https://aiesec.org/opportunity/global-talent/1331871'
$ws.Range("C2").Value = '[EXP] Digital Marketing PMO Intern'
$ws.Range("D2").Value = 'Maastricht, Netherlands'
$ws.Range("F2").Value = '16 applicants'
$ws.Range("G2").Value = '6 - 18 Months'
$ws.Range("H2").Value = 'DHL Group'

# ---------------------------------------------------------------------------
# 4) Row 3 - Financial Assistant Intern (Panamá)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "'1331873"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = 'https://aiesec.org/opportunity/global-talent/1331873'
$ws.Range("C3").Value = 'Financial Assistant Intern'
$ws.Range("D3").Value = 'Panamá, Provincia de Panamá, Panamá'
$ws.Range("F3").Value = '5 applicants'
$ws.Range("G3").Value = '6 - 18 Months'
$ws.Range("H3").Value = 'Nestlé'

# ---------------------------------------------------------------------------
# 5) Row 4 - Sales Representative (Istanbul, Türkiye)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "'1331867"
$ws.Range("A4").ClearFormats()
$ws.Range("B4").Value = 'https://aiesec.org/opportunity/global-talent/1331867'
$ws.Range("C4").Value = 'Sales Representative'
$ws.Range("D4").Value = 'Istanbul, İstanbul, Türkiye'
$ws.Range("F4").Value = '6 applicants'
$ws.Range("G4").Value = '3 - 6 Months'
$ws.Range("H4").Value = 'OPINNATE'

# ---------------------------------------------------------------------------
# 6) Row 5 - Procurement Intern (Dubai)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "'1331591"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value = 'https://aiesec.org/opportunity/global-talent/1331591'
$ws.Range("C5").Value = 'Procurement Intern'
$ws.Range("D5").Value = 'Dubai - United Arab Emirates'
$ws.Range("F5").Value = '32 applicants'
$ws.Range("G5").Value = '3 - 6 Months'
$ws.Range("H5").Value = 'Dubai Holding Group Services'

# ---------------------------------------------------------------------------
# 7) Row 6 - Order to Cash (Accounts Receivable) Intern (Dubai)
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "'1331590"
$ws.Range("A6").ClearFormats()
$ws.Range("B6").Value = 'https://aiesec.org/opportunity/global-talent/1331590'
$ws.Range("C6").Value = 'Order to Cash (Accounts Receivable) Intern'
$ws.Range("D6").Value = 'Dubai - United Arab Emirates'
$ws.Range("F6").Value = '16 applicants'
$ws.Range("G6").Value = '3 - 6 Months'
$ws.Range("H6").Value = 'Dubai Holding Group Services'

# ---------------------------------------------------------------------------
# 8) Row 7 - Master Data Management Intern (Dubai)
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "'1331468"
$ws.Range("A7").ClearFormats()
$ws.Range("B7").Value = 'https://aiesec.org/opportunity/global-talent/1331468'
$ws.Range("C7").Value = 'Master Data Management Intern'
$ws.Range("D7").Value = 'Dubai - United Arab Emirates'
$ws.Range("F7").Value = '15 applicants'
$ws.Range("G7").Value = '3 - 6 Months'
$ws.Range("H7").Value = 'Dubai Holding Group Services'

# ---------------------------------------------------------------------------
# 9) Row 8 - Accelerate Romania - Junior MEP Design Engineer (Sibiu)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "'1331119"
$ws.Range("A8").ClearFormats()
$ws.Range("B8").Value = 'https://aiesec.org/opportunity/global-talent/1331119'
$ws.Range("C8").Value = 'Accelerate Romania - Junior MEP Design Engineer'
$ws.Range("D8").Value = 'Sibiu, Romania'
$ws.Range("F8").Value = '0 applicants'
$ws.Range("H8").Value = 'GB Instaplan'

# ---------------------------------------------------------------------------
# 10) Row 9 - Digital Technology Application Services (ONLY EU) (Bruxelles)
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "'1330419"
$ws.Range("A9").ClearFormats()
$ws.Range("B9").Value = 'https://aiesec.org/opportunity/global-talent/1330419'
$ws.Range("C9").Value = 'Digital Technology Application Services (ONLY EU)'
$ws.Range("D9").Value = 'Bruxelles, Belgio'
$ws.Range("F9").Value = '46 applicants'
$ws.Range("G9").Value = '6 - 18 Months'
$ws.Range("H9").Value = 'UCB'

# ---------------------------------------------------------------------------
# 11) Row 10 - Resources humaines (Bouaké, Côte d'Ivoire)
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "'1323274"
$ws.Range("A10").ClearFormats()
$ws.Range("B10").Value = 'https://aiesec.org/opportunity/global-talent/1323274'
$ws.Range("C10").Value = 'Resources humaines'
$ws.Range("D10").Value = 'Bouaké, Côte d''Ivoire'
$ws.Range("F10").Value = '8 applicants'
$ws.Range("G10").Value = '9 - 12 Weeks'
$ws.Range("H10").Value = 'Radio Media Plus'

# ---------------------------------------------------------------------------
# 12) Column width tweaks (C, D, F, H got narrower to fit the new content).
#     XLSX stores a "characters" width that is 5/6 wider than the COM
#     ColumnWidth value, so subtract 5/6 to land on the exact target width.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 51.1666666666667
$ws.Columns("D").ColumnWidth = 37.1666666666667
$ws.Columns("F").ColumnWidth = 15.1666666666667
$ws.Columns("H").ColumnWidth = 30.1666666666667
